$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 41
$ws.Range("K6").Value = 123
$ws.Range("M6").Value = -11

$ws.Range("H32").Value = 3581.0667
$ws.Range("I32").Value = 2824.75
$ws.Range("K32").Value = 2824.75
$ws.Range("M32").Value = -2498.75

$ws.Range("H33").Value = 3704294.2
$ws.Range("J33").Value = 874
$ws.Range("L33").Value = 874
$ws.Range("N33").Value = -1332

$ws.Range("H39").Value = 209.94737
$ws.Range("I39").Value = 187.36363
$ws.Range("J39").Value = 241
$ws.Range("K39").Value = 562.0908899999999
$ws.Range("L39").Value = 723
$ws.Range("M39").Value = -266.0908899999999
$ws.Range("N39").Value = -1315

$ws.Range("H51").Value = 6373.467
$ws.Range("I51").Value = 4555.4443
$ws.Range("J51").Value = 7152.619
$ws.Range("K51").Value = 4555.4443
$ws.Range("L51").Value = 7152.619
$ws.Range("M51").Value = -4071.4443
$ws.Range("N51").Value = -8120.619

$ws.Range("H55").Value = 176.26666
$ws.Range("I55").Value = 175.5
$ws.Range("J55").Value = 177.14285
$ws.Range("K55").Value = 175.5
$ws.Range("L55").Value = 177.14285
$ws.Range("M55").Value = 38.5
$ws.Range("N55").Value = -605.14285

$ws.Range("H62").Value = 4982.147
$ws.Range("I62").Value = 4335.5
$ws.Range("J62").Value = 7999.8335
$ws.Range("K62").Value = 4335.5
$ws.Range("L62").Value = 7999.8335
$ws.Range("M62").Value = -3711.5
$ws.Range("N62").Value = -9247.833500000001

$ws.Range("H65").Value = 4982.147
$ws.Range("I65").Value = 4335.5
$ws.Range("J65").Value = 7999.8335
$ws.Range("K65").Value = 21677.5
$ws.Range("L65").Value = 39999.1675
$ws.Range("M65").Value = -18557.5
$ws.Range("N65").Value = -46239.1675

$ws.Range("H70").Value = 3879.6667
$ws.Range("I70").Value = 1792.3334
$ws.Range("J70").Value = 5967
$ws.Range("K70").Value = 5377.0002
$ws.Range("L70").Value = 17901
$ws.Range("M70").Value = -5107.0002
$ws.Range("N70").Value = -18441

$ws.Range("H73").Value = 3879.6667
$ws.Range("I73").Value = 1792.3334
$ws.Range("J73").Value = 5967
$ws.Range("K73").Value = 5377.0002
$ws.Range("L73").Value = 17901
$ws.Range("M73").Value = -4441.0002
$ws.Range("N73").Value = -19773

$ws.Range("H74").Value = 7437.393
$ws.Range("I74").Value = 5672.1665
$ws.Range("J74").Value = 7918.8184
$ws.Range("K74").Value = 5672.1665
$ws.Range("L74").Value = 7918.8184
$ws.Range("M74").Value = -4736.1665
$ws.Range("N74").Value = -9790.8184

$ws.Range("H76").Value = 5105.1113
$ws.Range("J76").Value = 6463.125
$ws.Range("L76").Value = 6463.125
$ws.Range("N76").Value = -7093.125

$ws.Range("H77").Value = 7437.393
$ws.Range("I77").Value = 5672.1665
$ws.Range("J77").Value = 7918.8184
$ws.Range("K77").Value = 28360.8325
$ws.Range("L77").Value = 39594.092
$ws.Range("M77").Value = -23680.8325
$ws.Range("N77").Value = -48954.092

$ws.Range("H79").Value = 5105.1113
$ws.Range("J79").Value = 6463.125
$ws.Range("L79").Value = 6463.125
$ws.Range("N79").Value = -8647.125

$ws.Range("H86").Value = 8451.362999999999
$ws.Range("I86").Value = 5324
$ws.Range("J86").Value = 9624.125
$ws.Range("K86").Value = 5324
$ws.Range("L86").Value = 9624.125
$ws.Range("M86").Value = -4201
$ws.Range("N86").Value = -11870.125

$ws.Range("H89").Value = 8451.362999999999
$ws.Range("I89").Value = 5324
$ws.Range("J89").Value = 9624.125
$ws.Range("K89").Value = 26620
$ws.Range("L89").Value = 48120.625
$ws.Range("M89").Value = -21004
$ws.Range("N89").Value = -59352.625

$ws.Range("H96").Value = 219.4375
$ws.Range("I96").Value = 230.27272
$ws.Range("J96").Value = 195.6
$ws.Range("K96").Value = 690.81816
$ws.Range("L96").Value = 586.8
$ws.Range("M96").Value = 682.18184
$ws.Range("N96").Value = -3332.8

$ws.Range("H98").Value = 1869.1818
$ws.Range("I98").Value = 1796.258
$ws.Range("J98").Value = 2999.5
$ws.Range("K98").Value = 1796.258
$ws.Range("L98").Value = 2999.5
$ws.Range("M98").Value = -298.258
$ws.Range("N98").Value = -5995.5

$ws.Range("H101").Value = 35714704
$ws.Range("I101").Value = 45454916
$ws.Range("K101").Value = 136364748
$ws.Range("M101").Value = -136363126

$ws.Range("H112").Value = 6781.1055
$ws.Range("J112").Value = 6781.1055
$ws.Range("L112").Value = 20343.3165
$ws.Range("N112").Value = -22559.3165

$ws.Range("H116").Value = 6532.75
$ws.Range("J116").Value = 7892
$ws.Range("L116").Value = 7892
$ws.Range("N116").Value = -14776

$ws.Range("H118").Value = 90909480
$ws.Range("I118").Value = 111111470
$ws.Range("J118").Value = 525
$ws.Range("K118").Value = 333334410
$ws.Range("L118").Value = 1575
$ws.Range("M118").Value = -333332753
$ws.Range("N118").Value = -4889

$ws.Range("H122").Value = 1869.1818
$ws.Range("I122").Value = 1796.258
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 5388.774
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -2938.774
$ws.Range("N122").Value = -13898.5

$ws.Range("H131").Value = 3310.1538
$ws.Range("I131").Value = 1126.8889
$ws.Range("J131").Value = 8222.5
$ws.Range("K131").Value = 3380.6667
$ws.Range("L131").Value = 24667.5
$ws.Range("M131").Value = 1659.3333
$ws.Range("N131").Value = -34747.5

$ws.Range("H137").Value = 76063.875
$ws.Range("I137").Value = 119920.13
$ws.Range("K137").Value = 359760.39
$ws.Range("M137").Value = -357210.39

$ws.Range("H138").Value = 3560.6938
$ws.Range("J138").Value = 3794.4358
$ws.Range("L138").Value = 11383.3074
$ws.Range("N138").Value = -21663.3074

$ws.Range("H139").Value = 120665.336
$ws.Range("J139").Value = 120998
$ws.Range("L139").Value = 120998
$ws.Range("N139").Value = -131278

$ws.Range("H141").Value = 20338.055
$ws.Range("I141").Value = 10739
$ws.Range("K141").Value = 32217
$ws.Range("M141").Value = -27037

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4113130.5
$ws.Range("J45").Value = 5515.8335
$ws.Range("L45").Value = 5515.8335
$ws.Range("N45").Value = -6269.8335

$ws.Range("H61").Value = 4291
$ws.Range("I61").Value = 4290.1763
$ws.Range("K61").Value = 4290.1763
$ws.Range("M61").Value = -4078.1763

$ws.Range("H74").Value = 47085.332
$ws.Range("I74").Value = 3398.4
$ws.Range("J74").Value = 156302.67
$ws.Range("K74").Value = 3398.4
$ws.Range("L74").Value = 156302.67
$ws.Range("M74").Value = -2524.4
$ws.Range("N74").Value = -158050.67

$ws.Range("H77").Value = 47085.332
$ws.Range("I77").Value = 3398.4
$ws.Range("J77").Value = 156302.67
$ws.Range("K77").Value = 16992
$ws.Range("L77").Value = 781513.3500000001
$ws.Range("M77").Value = -12624
$ws.Range("N77").Value = -790249.3500000001

$ws.Range("H88").Value = 2683.4285
$ws.Range("I88").Value = 2591
$ws.Range("K88").Value = 2591
$ws.Range("M88").Value = -2185

$ws.Range("H91").Value = 2683.4285
$ws.Range("I91").Value = 2591
$ws.Range("K91").Value = 2591
$ws.Range("M91").Value = -1187

$ws.Range("H102").Value = 2875867
$ws.Range("I102").Value = 3625068.5
$ws.Range("J102").Value = 3928.3333
$ws.Range("K102").Value = 3625068.5
$ws.Range("L102").Value = 3928.3333
$ws.Range("M102").Value = -3623446.5
$ws.Range("N102").Value = -7172.3333

$ws.Range("H103").Value = 58142
$ws.Range("J103").Value = 58142
$ws.Range("L103").Value = 58142
$ws.Range("N103").Value = -60486

$ws.Range("H136").Value = 4291
$ws.Range("I136").Value = 4290.1763
$ws.Range("K136").Value = 12870.5289
$ws.Range("M136").Value = -10320.5289

$ws.Range("H138").Value = 109996.5
$ws.Range("J138").Value = 109996.5
$ws.Range("L138").Value = 109996.5
$ws.Range("N138").Value = -120276.5

$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

$ws.Range("H140").Value = 118428
$ws.Range("J140").Value = 118428
$ws.Range("L140").Value = 118428
$ws.Range("N140").Value = -128788

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4168826
$ws.Range("I86").Value = 5557860.5
$ws.Range("J86").Value = 1721.5
$ws.Range("K86").Value = 5557860.5
$ws.Range("L86").Value = 1721.5
$ws.Range("M86").Value = -5556737.5
$ws.Range("N86").Value = -3967.5

$ws.Range("H89").Value = 4168826
$ws.Range("I89").Value = 5557860.5
$ws.Range("J89").Value = 1721.5
$ws.Range("K89").Value = 27789302.5
$ws.Range("L89").Value = 8607.5
$ws.Range("M89").Value = -27783686.5
$ws.Range("N89").Value = -19839.5

$ws.Range("H105").Value = 6946696.5
$ws.Range("J105").Value = 3472.25
$ws.Range("L105").Value = 3472.25
$ws.Range("N105").Value = -6966.25

$ws.Range("H132").Value = 89350
$ws.Range("J132").Value = 89350
$ws.Range("L132").Value = 89350
$ws.Range("N132").Value = -99470

$ws.Range("H134").Value = 8459.4
$ws.Range("I134").Value = 3364.6667
$ws.Range("J134").Value = 10642.857
$ws.Range("K134").Value = 10094.0001
$ws.Range("L134").Value = 31928.571
$ws.Range("M134").Value = -7559.000100000001
$ws.Range("N134").Value = -36998.571

$ws.Range("H135").Value = 75000
$ws.Range("J135").Value = 75000
$ws.Range("L135").Value = 75000
$ws.Range("N135").Value = -85140

$ws.Range("H140").Value = 88776.75
$ws.Range("J140").Value = 88776.75
$ws.Range("L140").Value = 88776.75
$ws.Range("N140").Value = -99136.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 948.4783
$ws.Range("I16").Value = 531.8182
$ws.Range("J16").Value = 1330.4166
$ws.Range("K16").Value = 531.8182
$ws.Range("L16").Value = 1330.4166
$ws.Range("M16").Value = -244.8182
$ws.Range("N16").Value = -1904.4166

$ws.Range("H31").Value = 21975.725
$ws.Range("I31").Value = 2321.6
$ws.Range("J31").Value = 30164.945
$ws.Range("K31").Value = 2321.6
$ws.Range("L31").Value = 30164.945
$ws.Range("M31").Value = -2026.6
$ws.Range("N31").Value = -30754.945

$ws.Range("H34").Value = 21975.725
$ws.Range("I34").Value = 2321.6
$ws.Range("J34").Value = 30164.945
$ws.Range("K34").Value = 2321.6
$ws.Range("L34").Value = 30164.945
$ws.Range("M34").Value = -2119.6
$ws.Range("N34").Value = -30568.945

$ws.Range("H62").Value = 3802
$ws.Range("I62").Value = 3802
$ws.Range("K62").Value = 3802
$ws.Range("M62").Value = -3178

$ws.Range("H65").Value = 3802
$ws.Range("I65").Value = 3802
$ws.Range("K65").Value = 19010
$ws.Range("M65").Value = -15890

$ws.Range("H105").Value = 1385.25
$ws.Range("J105").Value = 1943.375
$ws.Range("L105").Value = 1943.375
$ws.Range("N105").Value = -5437.375

$ws.Range("H107").Value = 3495.5833
$ws.Range("I107").Value = 3552.1
$ws.Range("K107").Value = 3552.1
$ws.Range("M107").Value = -1632.1

$ws.Range("H113").Value = 948.4783
$ws.Range("I113").Value = 531.8182
$ws.Range("J113").Value = 1330.4166
$ws.Range("K113").Value = 531.8182
$ws.Range("L113").Value = 1330.4166
$ws.Range("M113").Value = 1638.1818
$ws.Range("N113").Value = -5670.4166

$ws.Range("H120").Value = 100000
$ws.Range("J120").Value = 100000
$ws.Range("L120").Value = 100000
$ws.Range("N120").Value = -107258

$ws.Range("H122").Value = 2861.48
$ws.Range("I122").Value = 2860.5
$ws.Range("J122").Value = 2865.4
$ws.Range("K122").Value = 8581.5
$ws.Range("L122").Value = 8596.200000000001
$ws.Range("M122").Value = -6131.5
$ws.Range("N122").Value = -13496.2

$ws.Range("H141").Value = 990000
$ws.Range("J141").Value = 990000
$ws.Range("L141").Value = 990000
$ws.Range("N141").Value = -1000360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15337251
$ws.Range("I4").Value = 17795292
$ws.Range("K4").Value = 53385876
$ws.Range("M4").Value = -53385764

$ws.Range("H37").Value = 46162.5
$ws.Range("J37").Value = 46162.5
$ws.Range("L37").Value = 138487.5
$ws.Range("N37").Value = -138711.5

$ws.Range("H57").Value = 2965.3635
$ws.Range("I57").Value = 529.75
$ws.Range("K57").Value = 1589.25
$ws.Range("M57").Value = -1030.25

$ws.Range("H98").Value = 2669.3333
$ws.Range("J98").Value = 2669.3333
$ws.Range("L98").Value = 8007.999899999999
$ws.Range("N98").Value = -11003.9999

$ws.Range("H105").Value = 100
$ws.Range("I105").Value = 100
$ws.Range("K105").Value = 300
$ws.Range("M105").Value = 2321

$ws.Range("H121").Value = 2079.1428
$ws.Range("I121").Value = 518
$ws.Range("J121").Value = 3250
$ws.Range("K121").Value = 1554
$ws.Range("L121").Value = 9750
$ws.Range("M121").Value = -244
$ws.Range("N121").Value = -12370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 74916220
$ws.Range("I80").Value = 131099760
$ws.Range("K80").Value = 131099760
$ws.Range("M80").Value = -131098762

$ws.Range("H83").Value = 74916220
$ws.Range("I83").Value = 131099760
$ws.Range("K83").Value = 655498800
$ws.Range("M83").Value = -655493808

$ws.Range("H97").Value = 1488657
$ws.Range("I97").Value = 2976741
$ws.Range("J97").Value = 572.875
$ws.Range("K97").Value = 2976741
$ws.Range("L97").Value = 572.875
$ws.Range("M97").Value = -2976245
$ws.Range("N97").Value = -1564.875

$ws.Range("H122").Value = 475031
$ws.Range("I122").Value = 529387.5600000001
$ws.Range("K122").Value = 1588162.68
$ws.Range("M122").Value = -1585712.68

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3346666.2
$ws.Range("J2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("N2").Value = -20224

$ws.Range("H12").Value = 4499.875
$ws.Range("I12").Value = 2999.5
$ws.Range("K12").Value = 2999.5
$ws.Range("M12").Value = -2829.5

$ws.Range("H22").Value = 51022.945
$ws.Range("I22").Value = 69267.234
$ws.Range("J22").Value = 3587.8
$ws.Range("K22").Value = 69267.234
$ws.Range("L22").Value = 3587.8
$ws.Range("M22").Value = -68972.234
$ws.Range("N22").Value = -4177.8

$ws.Range("H27").Value = 51022.945
$ws.Range("I27").Value = 69267.234
$ws.Range("J27").Value = 3587.8
$ws.Range("K27").Value = 69267.234
$ws.Range("L27").Value = 3587.8
$ws.Range("M27").Value = -69160.234
$ws.Range("N27").Value = -3801.8

$ws.Range("H32").Value = 15013
$ws.Range("I32").Value = 15013
$ws.Range("K32").Value = 15013
$ws.Range("M32").Value = -14696

$ws.Range("H100").Value = 1933.1111
$ws.Range("I100").Value = 1567.9546
$ws.Range("J100").Value = 3539.8
$ws.Range("K100").Value = 1567.9546
$ws.Range("L100").Value = 3539.8
$ws.Range("M100").Value = -1026.9546
$ws.Range("N100").Value = -4621.8

$ws.Range("H122").Value = 4775.6665
$ws.Range("I122").Value = 2499.1333
$ws.Range("K122").Value = 7497.3999
$ws.Range("M122").Value = -5047.3999

$ws.Range("H132").Value = 5323.2417
$ws.Range("I132").Value = 5554.2114
$ws.Range("J132").Value = 4122.2
$ws.Range("K132").Value = 16662.6342
$ws.Range("L132").Value = 12366.6
$ws.Range("M132").Value = -14132.6342
$ws.Range("N132").Value = -17426.6

$ws.Range("H136").Value = 51154.395
$ws.Range("I136").Value = 70191.3
$ws.Range("K136").Value = 210573.9
$ws.Range("M136").Value = -208023.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = ""

$ws.Range("H113").Value = 1084.9584
$ws.Range("I113").Value = 541
$ws.Range("K113").Value = 1623
$ws.Range("M113").Value = 547

$ws.Range("H126").Value = 1890.2632
$ws.Range("I126").Value = 1982.1875
$ws.Range("J126").Value = 1400
$ws.Range("K126").Value = 5946.5625
$ws.Range("L126").Value = 4200
$ws.Range("M126").Value = -3476.5625
$ws.Range("N126").Value = -9140

$ws.Range("H141").Value = 67332.664
$ws.Range("J141").Value = 67332.664
$ws.Range("L141").Value = 67332.664
$ws.Range("N141").Value = -77692.664
